$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.285.21"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").Value = "2.631.01"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'604.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'151.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("D10").Value = "'5.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.06%  "

$ws.Range("E11").Value = "  +6.70%  "

$ws.Range("D13").Value = "'27.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "3.103.92"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").Value = "64.137.99"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("E16").Value = "  +4.33%  "

$ws.Range("D17").Value = "2.601.90"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "'12.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.12%  "

$ws.Range("E19").Value = "  +4.40%  "

$ws.Range("D20").Value = "'350.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.28%  "

$ws.Range("D21").Value = "'7.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'5.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.86%  "

$ws.Range("D24").Value = "'66.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("E25").Value = "  +15.06%  "

$ws.Range("E26").Value = "  +5.45%  "

$ws.Range("D27").Value = "'9.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.73%  "

$ws.Range("E28").Value = "  +2.71%  "

$ws.Range("D29").Value = "'8.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.20%  "

$ws.Range("D30").Value = "'545.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").Value = "'2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").Value = "0.0₃0859"
$ws.Range("E33").Value = "  +7.44%  "

$ws.Range("D34").Value = "'1.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'5.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").Value = "'167.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").Value = "'2.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.16%  "

$ws.Range("D38").Value = "'0.411"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.41%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'19.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("D41").Value = "'172.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.86%  "

$ws.Range("D43").Value = "'40.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("D45").Value = "'0.0587"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.17%  "

$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").Value = "'0.629"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("E48").Value = "  +15.08%  "

$ws.Range("D49").Value = "'0.0247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.79%  "

$ws.Range("D50").Value = "'0.0968"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").Value = "'19.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.07%  "
